# Updated remaining queries for C3DC
# Fix the JOIN conditions in all SQL queries on Sheet1: the "id" columns
# are renamed to their fully-qualified equivalents (study_id / participant_id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Text
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $rng.Value = $text
}

# Widen column C to fit the longer query text (no longer an exact "best fit").
# Excel's ColumnWidth COM property and the stored <col width> XML attribute are
# offset by ~5/6 of a character (MDW rounding), so compensate to land on 67.5.
$ws.Columns.Item(3).ColumnWidth = 67.5 - (5/6)
